# The deck's slide-master design is switched from the "Integral" (Red
# Violet) theme colours over to the plain default "Office Theme" colour
# palette (ppt/theme/theme1.xml's <a:clrScheme> red-violet values ->
# the stock Office values: dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint's ColorScheme object exposes exactly those twelve theme
# slots (in this fixed order) through Colors(1..12), each an RGBColor
# whose .RGB is the classic VBA RGB() encoding: R + G*256 + B*65536.

function Get-RgbValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette (target), in ColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = Get-RgbValue $officeThemeColors[$i - 1]
}
